$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The shared strings used as status markers / labels are being refactored:
#   "⬛" (black square) -> "📘" (blue book)
#   "🟥" (red square)   -> "📕" (red book)
#   "🟩" (green square) -> "📗" (green book)
#   "🟧" (orange square)-> "📙" (orange book)
#   "noir" (black)      -> "bleu" (blue)
# Every cell in column A / B that previously held one of these values gets
# updated to its replacement, keeping everything else untouched.

$colAOld = "⬛"
$colANew = "📘"

$colAOldCells = @("A2", "A4", "A5", "A6", "A9", "A10", "A11", "A12", "A13", "A14")
foreach ($cell in $colAOldCells) {
    $ws.Range($cell).Value = $colANew
}

$ws.Range("A3").Value = "📕"
$ws.Range("A7").Value = "📗"
$ws.Range("A8").Value = "📙"

$colBOldCells = @("B2", "B4", "B5", "B6", "B9", "B10", "B11", "B12", "B13", "B14")
foreach ($cell in $colBOldCells) {
    $ws.Range($cell).Value = "bleu"
}
